# Insert a new data row at row 41 (pushing the existing rows 41..124 down to
# 42..125) and populate it with a new "Cilantro" price observation.
#
# Columns are: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificación.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 41:124 down to 42:125, leaving row 41 free for the new record.
$ws.Rows.Item(41).Insert()

# Match the date-cell number format used by the rest of column D.
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat

$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 45152
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112040
$ws.Range("G41").Value = "Cilantro"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 250
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 3500
$ws.Range("M41").Value = 3250
$ws.Range("N41").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 1625
$ws.Range("Q41").Value = 2
$ws.Range("R41").Value = "Hortaliza"
